# Daily attendance processing - clear "Recorded By" (column G) values
# and shrink column G width now that the names are no longer shown.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Row 1 holds the column headers ("Recorded By") - leave it untouched and
# only clear the data rows below it.
$startRow = 2
if ($firstRow -gt $startRow) {
    $startRow = $firstRow
}

for ($r = $startRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -ne "") {
        $cell.Value = ""
    }
}

# Shrink the "Recorded By" column now that it's empty (was 50 chars wide,
# now 13 chars wide).
$ws.Columns.Item(7).ColumnWidth = 12.17
